$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for columns G (Return_with_prediction) and H (return_pct_change)
# for data rows 2-85, plus one updated value in column I (mean_return_pct_change) on row 2.
$newValues = @(
    @{ Row = 2; G = 0.0687161128519025; H = 6.864325880224633 }
    @{ Row = 3; G = 0.06962990972764119; H = 23.88111868055159 }
    @{ Row = 4; G = -0.01781769405652903; H = 12.27826956698774 }
    @{ Row = 5; G = -0.01279179220849823; H = -11.90116223971244 }
    @{ Row = 6; G = -0.01072865862316981; H = 4.153273335132798 }
    @{ Row = 7; G = -0.001035909746857211; H = 81.7097805263429 }
    @{ Row = 8; G = -0.01166555266453533; H = -105.0486630511336 }
    @{ Row = 9; G = -0.01080026839265989; H = -96.61352393937428 }
    @{ Row = 10; G = -0.06945140951564707; H = -10.20227840086546 }
    @{ Row = 11; G = -0.07088649740807156; H = -10.62524473530545 }
    @{ Row = 12; G = -0.3801062574092129; H = 3.673851621203487 }
    @{ Row = 13; G = -0.4033861550113001; H = -2.87334944754482 }
    @{ Row = 14; G = -0.01791969331787071; H = -118.6858316637777 }
    @{ Row = 15; G = -0.02517251015128741; H = 44.45444422926664 }
    @{ Row = 16; G = 0.137479837345001; H = 0.5339596934981469 }
    @{ Row = 17; G = 0.1463664972265752; H = 4.944308554069225 }
    @{ Row = 18; G = 0.118676273737338; H = 0.8567020210108911 }
    @{ Row = 19; G = 0.1247092161422848; H = -3.077384351831464 }
    @{ Row = 20; G = 0.08562269976166079; H = -3.507010553650305 }
    @{ Row = 21; G = 0.08551956517455996; H = -1.804388096273384 }
    @{ Row = 22; G = -0.1012253982017692; H = -8.281949407746087 }
    @{ Row = 23; G = -0.09776948567779378; H = 3.623176883872968 }
    @{ Row = 24; G = 0.1579489941306083; H = -1.950491434747088 }
    @{ Row = 25; G = 0.169370734164759; H = -0.7180706012560135 }
    @{ Row = 26; G = 0.09968957830647295; H = 9.967351668985296 }
    @{ Row = 27; G = 0.09043425337695868; H = 5.202903955464625 }
    @{ Row = 28; G = -0.1304852900002428; H = 5.198433360166433 }
    @{ Row = 29; G = -0.1402205903515006; H = -0.3010197344406161 }
    @{ Row = 30; G = 0.0477613497369185; H = -8.178222428350121 }
    @{ Row = 31; G = 0.04858110923370425; H = 10.86666155482429 }
    @{ Row = 32; G = 0.1162055733462721; H = 6.900984689387695 }
    @{ Row = 33; G = 0.1152263702669514; H = -7.146789490141786 }
    @{ Row = 34; G = -0.01328841399165122; H = 14.90728245733816 }
    @{ Row = 35; G = -0.0132704151696055; H = 20.71347014645429 }
    @{ Row = 36; G = 0.03530580708966877; H = -3.973960420367229 }
    @{ Row = 37; G = 0.03435633251965751; H = -3.725071564648125 }
    @{ Row = 38; G = 0.09931858199055098; H = -0.9806834708436681 }
    @{ Row = 39; G = 0.09458493476401787; H = -2.890225638655779 }
    @{ Row = 40; G = 0.0329659758490939; H = -2.142163511107762 }
    @{ Row = 41; G = 0.03147227944696478; H = -2.316806722155458 }
    @{ Row = 42; G = 0.1189971445072512; H = -1.578329418707107 }
    @{ Row = 43; G = 0.1269318328390656; H = -0.6680059347719242 }
    @{ Row = 44; G = 0.04012999885591129; H = 1.16813435378751 }
    @{ Row = 45; G = 0.03107227238396221; H = -0.2890685034152298 }
    @{ Row = 46; G = 0.06192739713121615; H = 9.378225187172283 }
    @{ Row = 47; G = 0.06423241472782022; H = 9.480501736196121 }
    @{ Row = 48; G = 0.04320505156588467; H = -12.28704854835148 }
    @{ Row = 49; G = 0.05107741030319914; H = 12.09114637261196 }
    @{ Row = 50; G = 0.02537628235974115; H = -4.191577123390303 }
    @{ Row = 51; G = 0.02111242018220849; H = -24.64120749155786 }
    @{ Row = 52; G = -0.08059614848701571; H = 7.279369045327234 }
    @{ Row = 53; G = -0.08056311380526139; H = -0.4324710817105894 }
    @{ Row = 54; G = 0.05074865819560606; H = 1.443141843197838 }
    @{ Row = 55; G = 0.05378065964043954; H = -4.449488515150281 }
    @{ Row = 56; G = 0.04758644232247539; H = -3.733323728292866 }
    @{ Row = 57; G = 0.0473998097331657; H = 24.78895683378744 }
    @{ Row = 58; G = 0.05702397961730128; H = -1.024182409867418 }
    @{ Row = 59; G = 0.05450743124128897; H = -4.395250990200509 }
    @{ Row = 60; G = 0.03064826071491992; H = 11.61056537647907 }
    @{ Row = 61; G = 0.03199557913181796; H = 19.84339080852636 }
    @{ Row = 62; G = 0.06446805011198303; H = 3.236042333559309 }
    @{ Row = 63; G = 0.06781019502107434; H = 6.148029984641779 }
    @{ Row = 64; G = 0.02858970992060214; H = 3.064394357718167 }
    @{ Row = 65; G = 0.03685844230668094; H = 4.039800765888352 }
    @{ Row = 66; G = 0.07671251759413329; H = -1.249515734495227 }
    @{ Row = 67; G = 0.0845408478024388; H = 7.188987899267529 }
    @{ Row = 68; G = -0.02498525822933346; H = -14.89750423500255 }
    @{ Row = 69; G = -0.0260478217823271; H = -36.08169343397281 }
    @{ Row = 70; G = 0.07086713855102622; H = -1.554088864819031 }
    @{ Row = 71; G = 0.07987693348424307; H = 0.5740237726624602 }
    @{ Row = 72; G = -0.1440877469593357; H = 6.218200480364896 }
    @{ Row = 73; G = -0.144243303121051; H = 5.769864250332136 }
    @{ Row = 74; G = 0.1517473851888521; H = 0.8829932434955047 }
    @{ Row = 75; G = 0.1571684401280606; H = 4.467421145361128 }
    @{ Row = 76; G = -0.008067577993101892; H = -678.3768862086888 }
    @{ Row = 77; G = -0.008625555004576185; H = -290.6699907092726 }
    @{ Row = 78; G = 0.08539625475252073; H = -5.071236537967876 }
    @{ Row = 79; G = 0.09402122229861852; H = -2.97224134535508 }
    @{ Row = 80; G = -0.2307425332813548; H = -6.611349052602579 }
    @{ Row = 81; G = -0.2028755601880298; H = 4.804156772934314 }
    @{ Row = 82; G = 0.1740152756391153; H = 3.816785845426194 }
    @{ Row = 83; G = 0.1809372884403179; H = 2.78838040696195 }
    @{ Row = 84; G = 0.1099987507459851; H = 3.645590103073284 }
    @{ Row = 85; G = 0.1032876397970545; H = -1.215743700015579 }
)

foreach ($entry in $newValues) {
    $ws.Cells.Item($entry.Row, 7).Value = $entry.G   # column G
    $ws.Cells.Item($entry.Row, 8).Value = $entry.H   # column H
}

# Column I (mean_return_pct_change) only changes on row 2
$ws.Cells.Item(2, 9).Value = -12.96281062771407

